$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells C1 and D1, copying the style of B1 (bold header style)
$ws.Range("C1").Value = "valor_iptu"
$ws.Range("D1").Value = "valor_condominio"
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Update row 2 data
$ws.Range("A2").Value = "dazin"
$ws.Range("B2").Value = 1500
$ws.Range("C2").Value = 300
$ws.Range("D2").Value = 300

# Delete row 3 entirely (Joao / 550)
$ws.Rows(3).Delete()
